$d = $word.ActiveDocument

# The first paragraph currently reads:
#   "This is a Microsoft word document."
# It should become (as three additional, separately-represented runs):
#   "This is a Microsoft word document." + " (" + "Changed main" + ")"
# i.e. "This is a Microsoft word document. (Changed main)"

$para = $d.Paragraphs.Item(1)
$r = $para.Range
$r.End = $r.End - 1   # exclude the paragraph mark; keep only the run text

# Use InsertXML (Range.InsertXML / WordOpenXML) so the new text lands in its
# own separate <w:r> runs instead of being merged into the existing run's
# text, matching the four-run structure produced by the original edit.
$xml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t>This is a Microsoft word document.</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:t>Changed main</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$r.InsertXML($xml)
